$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns where values must stay literal text (dates as text, "NN%" as text, result codes)
$textCols = @(1,2,3,17,18,19)
foreach ($col in $textCols) {
    $ws.Range($ws.Cells.Item(163, $col), $ws.Cells.Item(179, $col)).NumberFormat = "@"
}

# Row 163
$ws.Cells.Item(163, 1).Value = '2025-07-18'
$ws.Cells.Item(163, 2).Value = 'remo'
$ws.Cells.Item(163, 3).Value = 'Novorizontino'
$ws.Cells.Item(163, 4).Value = 1
$ws.Cells.Item(163, 5).Value = 1
$ws.Cells.Item(163, 6).Value = 1353423
$ws.Cells.Item(163, 7).Value = 5
$ws.Cells.Item(163, 8).Value = 7
$ws.Cells.Item(163, 9).Value = 3
$ws.Cells.Item(163, 10).Value = 2
$ws.Cells.Item(163, 11).Value = 0
$ws.Cells.Item(163, 12).Value = 0
$ws.Cells.Item(163, 13).Value = 0
$ws.Cells.Item(163, 14).Value = 0
$ws.Cells.Item(163, 15).Value = 1
$ws.Cells.Item(163, 16).Value = 1
$ws.Cells.Item(163, 17).Value = '52%'
$ws.Cells.Item(163, 18).Value = '48%'
$ws.Cells.Item(163, 19).Value = 'E'

# Row 164
$ws.Cells.Item(164, 1).Value = '2025-07-18'
$ws.Cells.Item(164, 2).Value = 'Atletico Goianiense'
$ws.Cells.Item(164, 3).Value = 'Criciuma'
$ws.Cells.Item(164, 4).Value = 0
$ws.Cells.Item(164, 5).Value = 1
$ws.Cells.Item(164, 6).Value = 1353422
$ws.Cells.Item(164, 7).Value = 4
$ws.Cells.Item(164, 8).Value = 2
$ws.Cells.Item(164, 9).Value = 0
$ws.Cells.Item(164, 10).Value = 2
$ws.Cells.Item(164, 11).Value = 0
$ws.Cells.Item(164, 12).Value = 0
$ws.Cells.Item(164, 13).Value = 0
$ws.Cells.Item(164, 14).Value = 0
$ws.Cells.Item(164, 15).Value = 0
$ws.Cells.Item(164, 16).Value = 1
$ws.Cells.Item(164, 17).Value = '56%'
$ws.Cells.Item(164, 18).Value = '44%'
$ws.Cells.Item(164, 19).Value = 'V'

# Row 165
$ws.Cells.Item(165, 1).Value = '2025-07-19'
$ws.Cells.Item(165, 2).Value = 'Ferroviária'
$ws.Cells.Item(165, 3).Value = 'Athletic Club'
$ws.Cells.Item(165, 4).Value = 1
$ws.Cells.Item(165, 5).Value = 2
$ws.Cells.Item(165, 6).Value = 1353428
$ws.Cells.Item(165, 7).Value = 2
$ws.Cells.Item(165, 8).Value = 2
$ws.Cells.Item(165, 9).Value = 7
$ws.Cells.Item(165, 10).Value = 2
$ws.Cells.Item(165, 11).Value = 4
$ws.Cells.Item(165, 12).Value = 0
$ws.Cells.Item(165, 13).Value = 0
$ws.Cells.Item(165, 14).Value = 0
$ws.Cells.Item(165, 15).Value = 1
$ws.Cells.Item(165, 16).Value = 2
$ws.Cells.Item(165, 17).Value = '40%'
$ws.Cells.Item(165, 18).Value = '60%'
$ws.Cells.Item(165, 19).Value = 'V'

# Row 166
$ws.Cells.Item(166, 1).Value = '2025-07-19'
$ws.Cells.Item(166, 2).Value = 'Goias'
$ws.Cells.Item(166, 3).Value = 'Cuiaba'
$ws.Cells.Item(166, 4).Value = 3
$ws.Cells.Item(166, 5).Value = 1
$ws.Cells.Item(166, 6).Value = 1353421
$ws.Cells.Item(166, 7).Value = 2
$ws.Cells.Item(166, 8).Value = 5
$ws.Cells.Item(166, 9).Value = 4
$ws.Cells.Item(166, 10).Value = 3
$ws.Cells.Item(166, 11).Value = 0
$ws.Cells.Item(166, 12).Value = 1
$ws.Cells.Item(166, 13).Value = 0
$ws.Cells.Item(166, 14).Value = 0
$ws.Cells.Item(166, 15).Value = 3
$ws.Cells.Item(166, 16).Value = 1
$ws.Cells.Item(166, 17).Value = '65%'
$ws.Cells.Item(166, 18).Value = '35%'
$ws.Cells.Item(166, 19).Value = 'L'

# Row 167
$ws.Cells.Item(167, 1).Value = '2025-07-19'
$ws.Cells.Item(167, 2).Value = 'Avai'
$ws.Cells.Item(167, 3).Value = 'Vila Nova'
$ws.Cells.Item(167, 4).Value = 1
$ws.Cells.Item(167, 5).Value = 1
$ws.Cells.Item(167, 6).Value = 1353424
$ws.Cells.Item(167, 7).Value = 10
$ws.Cells.Item(167, 8).Value = 6
$ws.Cells.Item(167, 9).Value = 0
$ws.Cells.Item(167, 10).Value = 3
$ws.Cells.Item(167, 11).Value = 0
$ws.Cells.Item(167, 12).Value = 0
$ws.Cells.Item(167, 13).Value = 0
$ws.Cells.Item(167, 14).Value = 0
$ws.Cells.Item(167, 15).Value = 1
$ws.Cells.Item(167, 16).Value = 1
$ws.Cells.Item(167, 17).Value = '47%'
$ws.Cells.Item(167, 18).Value = '53%'
$ws.Cells.Item(167, 19).Value = 'E'

# Row 168
$ws.Cells.Item(168, 1).Value = '2025-07-19'
$ws.Cells.Item(168, 2).Value = 'Coritiba'
$ws.Cells.Item(168, 3).Value = 'Paysandu'
$ws.Cells.Item(168, 4).Value = 2
$ws.Cells.Item(168, 5).Value = 5
$ws.Cells.Item(168, 6).Value = 1353420
$ws.Cells.Item(168, 7).Value = 10
$ws.Cells.Item(168, 8).Value = 2
$ws.Cells.Item(168, 9).Value = 4
$ws.Cells.Item(168, 10).Value = 2
$ws.Cells.Item(168, 11).Value = 0
$ws.Cells.Item(168, 12).Value = 0
$ws.Cells.Item(168, 13).Value = 0
$ws.Cells.Item(168, 14).Value = 0
$ws.Cells.Item(168, 15).Value = 2
$ws.Cells.Item(168, 16).Value = 5
$ws.Cells.Item(168, 17).Value = '67%'
$ws.Cells.Item(168, 18).Value = '33%'
$ws.Cells.Item(168, 19).Value = 'V'

# Row 169
$ws.Cells.Item(169, 1).Value = '2025-07-19'
$ws.Cells.Item(169, 2).Value = 'Volta Redonda'
$ws.Cells.Item(169, 3).Value = 'Atletico Paranaense'
$ws.Cells.Item(169, 4).Value = 3
$ws.Cells.Item(169, 5).Value = 2
$ws.Cells.Item(169, 6).Value = 1353429
$ws.Cells.Item(169, 7).Value = 3
$ws.Cells.Item(169, 8).Value = 4
$ws.Cells.Item(169, 9).Value = 3
$ws.Cells.Item(169, 10).Value = 5
$ws.Cells.Item(169, 11).Value = 0
$ws.Cells.Item(169, 12).Value = 1
$ws.Cells.Item(169, 13).Value = 0
$ws.Cells.Item(169, 14).Value = 0
$ws.Cells.Item(169, 15).Value = 3
$ws.Cells.Item(169, 16).Value = 2
$ws.Cells.Item(169, 17).Value = '42%'
$ws.Cells.Item(169, 18).Value = '58%'
$ws.Cells.Item(169, 19).Value = 'L'

# Row 170
$ws.Cells.Item(170, 1).Value = '2025-07-20'
$ws.Cells.Item(170, 2).Value = 'Amazonas'
$ws.Cells.Item(170, 3).Value = 'Botafogo SP'
$ws.Cells.Item(170, 4).Value = 3
$ws.Cells.Item(170, 5).Value = 0
$ws.Cells.Item(170, 6).Value = 1353426
$ws.Cells.Item(170, 7).Value = 1
$ws.Cells.Item(170, 8).Value = 7
$ws.Cells.Item(170, 9).Value = 2
$ws.Cells.Item(170, 10).Value = 1
$ws.Cells.Item(170, 11).Value = 0
$ws.Cells.Item(170, 12).Value = 0
$ws.Cells.Item(170, 13).Value = 0
$ws.Cells.Item(170, 14).Value = 0
$ws.Cells.Item(170, 15).Value = 3
$ws.Cells.Item(170, 16).Value = 0
$ws.Cells.Item(170, 17).Value = '42%'
$ws.Cells.Item(170, 18).Value = '58%'
$ws.Cells.Item(170, 19).Value = 'L'

# Row 171
$ws.Cells.Item(171, 1).Value = '2025-07-20'
$ws.Cells.Item(171, 2).Value = 'America Mineiro'
$ws.Cells.Item(171, 3).Value = 'Chapecoense-sc'
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 5).Value = 1
$ws.Cells.Item(171, 6).Value = 1353427
$ws.Cells.Item(171, 7).Value = 7
$ws.Cells.Item(171, 8).Value = 3
$ws.Cells.Item(171, 9).Value = 1
$ws.Cells.Item(171, 10).Value = 3
$ws.Cells.Item(171, 11).Value = 0
$ws.Cells.Item(171, 12).Value = 0
$ws.Cells.Item(171, 13).Value = 0
$ws.Cells.Item(171, 14).Value = 0
$ws.Cells.Item(171, 15).Value = 0
$ws.Cells.Item(171, 16).Value = 1
$ws.Cells.Item(171, 17).Value = '66%'
$ws.Cells.Item(171, 18).Value = '34%'
$ws.Cells.Item(171, 19).Value = 'V'

# Row 172
$ws.Cells.Item(172, 1).Value = '2025-07-22'
$ws.Cells.Item(172, 2).Value = 'Operario-PR'
$ws.Cells.Item(172, 3).Value = 'Atletico Goianiense'
$ws.Cells.Item(172, 4).Value = 3
$ws.Cells.Item(172, 5).Value = 0
$ws.Cells.Item(172, 6).Value = 1353435
$ws.Cells.Item(172, 7).Value = 3
$ws.Cells.Item(172, 8).Value = 2
$ws.Cells.Item(172, 9).Value = 3
$ws.Cells.Item(172, 10).Value = 6
$ws.Cells.Item(172, 11).Value = 0
$ws.Cells.Item(172, 12).Value = 2
$ws.Cells.Item(172, 13).Value = 0
$ws.Cells.Item(172, 14).Value = 0
$ws.Cells.Item(172, 15).Value = 3
$ws.Cells.Item(172, 16).Value = 0
$ws.Cells.Item(172, 17).Value = '54%'
$ws.Cells.Item(172, 18).Value = '46%'
$ws.Cells.Item(172, 19).Value = 'L'

# Row 173
$ws.Cells.Item(173, 1).Value = '2025-07-22'
$ws.Cells.Item(173, 2).Value = 'Atletico Paranaense'
$ws.Cells.Item(173, 3).Value = 'Ferroviária'
$ws.Cells.Item(173, 4).Value = 1
$ws.Cells.Item(173, 5).Value = 1
$ws.Cells.Item(173, 6).Value = 1353430
$ws.Cells.Item(173, 7).Value = 11
$ws.Cells.Item(173, 8).Value = 1
$ws.Cells.Item(173, 9).Value = 1
$ws.Cells.Item(173, 10).Value = 4
$ws.Cells.Item(173, 11).Value = 0
$ws.Cells.Item(173, 12).Value = 0
$ws.Cells.Item(173, 13).Value = 0
$ws.Cells.Item(173, 14).Value = 0
$ws.Cells.Item(173, 15).Value = 1
$ws.Cells.Item(173, 16).Value = 1
$ws.Cells.Item(173, 17).Value = '59%'
$ws.Cells.Item(173, 18).Value = '41%'
$ws.Cells.Item(173, 19).Value = 'E'

# Row 174
$ws.Cells.Item(174, 1).Value = '2025-07-23'
$ws.Cells.Item(174, 2).Value = 'Athletic Club'
$ws.Cells.Item(174, 3).Value = 'Coritiba'
$ws.Cells.Item(174, 4).Value = 1
$ws.Cells.Item(174, 5).Value = 1
$ws.Cells.Item(174, 6).Value = 1353437
$ws.Cells.Item(174, 7).Value = 5
$ws.Cells.Item(174, 8).Value = 3
$ws.Cells.Item(174, 9).Value = 2
$ws.Cells.Item(174, 10).Value = 3
$ws.Cells.Item(174, 11).Value = 0
$ws.Cells.Item(174, 12).Value = 0
$ws.Cells.Item(174, 13).Value = 0
$ws.Cells.Item(174, 14).Value = 0
$ws.Cells.Item(174, 15).Value = 1
$ws.Cells.Item(174, 16).Value = 1
$ws.Cells.Item(174, 17).Value = '42%'
$ws.Cells.Item(174, 18).Value = '58%'
$ws.Cells.Item(174, 19).Value = 'E'

# Row 175
$ws.Cells.Item(175, 1).Value = '2025-07-23'
$ws.Cells.Item(175, 2).Value = 'Vila Nova'
$ws.Cells.Item(175, 3).Value = 'CRB'
$ws.Cells.Item(175, 4).Value = 2
$ws.Cells.Item(175, 5).Value = 0
$ws.Cells.Item(175, 6).Value = 1353431
$ws.Cells.Item(175, 7).Value = 5
$ws.Cells.Item(175, 8).Value = 7
$ws.Cells.Item(175, 9).Value = 2
$ws.Cells.Item(175, 10).Value = 2
$ws.Cells.Item(175, 11).Value = 0
$ws.Cells.Item(175, 12).Value = 0
$ws.Cells.Item(175, 13).Value = 0
$ws.Cells.Item(175, 14).Value = 0
$ws.Cells.Item(175, 15).Value = 2
$ws.Cells.Item(175, 16).Value = 0
$ws.Cells.Item(175, 17).Value = '28%'
$ws.Cells.Item(175, 18).Value = '72%'
$ws.Cells.Item(175, 19).Value = 'L'

# Row 176
$ws.Cells.Item(176, 1).Value = '2025-07-23'
$ws.Cells.Item(176, 2).Value = 'Chapecoense-sc'
$ws.Cells.Item(176, 3).Value = 'Volta Redonda'
$ws.Cells.Item(176, 4).Value = 4
$ws.Cells.Item(176, 5).Value = 2
$ws.Cells.Item(176, 6).Value = 1353434
$ws.Cells.Item(176, 7).Value = 7
$ws.Cells.Item(176, 8).Value = 5
$ws.Cells.Item(176, 9).Value = 2
$ws.Cells.Item(176, 10).Value = 3
$ws.Cells.Item(176, 11).Value = 0
$ws.Cells.Item(176, 12).Value = 0
$ws.Cells.Item(176, 13).Value = 0
$ws.Cells.Item(176, 14).Value = 0
$ws.Cells.Item(176, 15).Value = 4
$ws.Cells.Item(176, 16).Value = 2
$ws.Cells.Item(176, 17).Value = '49%'
$ws.Cells.Item(176, 18).Value = '51%'
$ws.Cells.Item(176, 19).Value = 'L'

# Row 177
$ws.Cells.Item(177, 1).Value = '2025-07-23'
$ws.Cells.Item(177, 2).Value = 'Novorizontino'
$ws.Cells.Item(177, 3).Value = 'Goias'
$ws.Cells.Item(177, 4).Value = 1
$ws.Cells.Item(177, 5).Value = 0
$ws.Cells.Item(177, 6).Value = 1353438
$ws.Cells.Item(177, 7).Value = 3
$ws.Cells.Item(177, 8).Value = 8
$ws.Cells.Item(177, 9).Value = 4
$ws.Cells.Item(177, 10).Value = 2
$ws.Cells.Item(177, 11).Value = 0
$ws.Cells.Item(177, 12).Value = 0
$ws.Cells.Item(177, 13).Value = 0
$ws.Cells.Item(177, 14).Value = 0
$ws.Cells.Item(177, 15).Value = 1
$ws.Cells.Item(177, 16).Value = 0
$ws.Cells.Item(177, 17).Value = '47%'
$ws.Cells.Item(177, 18).Value = '53%'
$ws.Cells.Item(177, 19).Value = 'L'

# Row 178
$ws.Cells.Item(178, 1).Value = '2025-07-24'
$ws.Cells.Item(178, 2).Value = 'Cuiaba'
$ws.Cells.Item(178, 3).Value = 'America Mineiro'
$ws.Cells.Item(178, 4).Value = 3
$ws.Cells.Item(178, 5).Value = 1
$ws.Cells.Item(178, 6).Value = 1353432
$ws.Cells.Item(178, 7).Value = 2
$ws.Cells.Item(178, 8).Value = 9
$ws.Cells.Item(178, 9).Value = 3
$ws.Cells.Item(178, 10).Value = 2
$ws.Cells.Item(178, 11).Value = 0
$ws.Cells.Item(178, 12).Value = 1
$ws.Cells.Item(178, 13).Value = 0
$ws.Cells.Item(178, 14).Value = 0
$ws.Cells.Item(178, 15).Value = 3
$ws.Cells.Item(178, 16).Value = 1
$ws.Cells.Item(178, 17).Value = '34%'
$ws.Cells.Item(178, 18).Value = '66%'
$ws.Cells.Item(178, 19).Value = 'L'

# Row 179
$ws.Cells.Item(179, 1).Value = '2025-07-24'
$ws.Cells.Item(179, 2).Value = 'Amazonas'
$ws.Cells.Item(179, 3).Value = 'Paysandu'
$ws.Cells.Item(179, 4).Value = 1
$ws.Cells.Item(179, 5).Value = 1
$ws.Cells.Item(179, 6).Value = 1353436
$ws.Cells.Item(179, 7).Value = 13
$ws.Cells.Item(179, 8).Value = 3
$ws.Cells.Item(179, 9).Value = 3
$ws.Cells.Item(179, 10).Value = 0
$ws.Cells.Item(179, 11).Value = 0
$ws.Cells.Item(179, 12).Value = 1
$ws.Cells.Item(179, 13).Value = 0
$ws.Cells.Item(179, 14).Value = 0
$ws.Cells.Item(179, 15).Value = 1
$ws.Cells.Item(179, 16).Value = 1
$ws.Cells.Item(179, 17).Value = '68%'
$ws.Cells.Item(179, 18).Value = '32%'
$ws.Cells.Item(179, 19).Value = 'E'
